$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lista de activos")

# --- "Costo" (J) column: align number formats for the header/data-entry styles ---
# Header cell J3 (fill+border family) gets the "$ #,##0.0" style
$ws.Range("J3").NumberFormat = """$""#,##0.0;-""$""#,##0.0"
# Data-entry cell J4 (border+protected family) gets the red-negative "$ #,###.##000" style
$ws.Range("J4").NumberFormat = """$""#,###.##000_);[Red]\(""$""#,###.##000\)"
$ws.Range("J4").Value = 5453535345.65756

# --- New sample values below, using the column's existing default format ---
$ws.Range("J5").Value = [double]"3.45343453434345e18"
$ws.Range("J6").Value = [double]"4.35345345345345e21"
$ws.Range("J7").Value = 34543543.345434

# --- Row 1: the redundant (no-op) D1:F1 cell entries are cleared away ---
$ws.Range("D1:F1").Clear()

# --- Update the active selection to reflect where the user ended up ---
$ws.Range("J5").Select()

# --- "Categorías" sheet: make B11:C11 match the plain thin-border style used by the rows above (B4:C10) ---
$ws2 = $wb.Worksheets.Item("Categorías")
$rng = $ws2.Range("B11:C11")
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2
